# Insert a new data row at row 248 (this pushes the existing rows 248-344
# down to 249-345, and extends the used range to A1:R345), then populate
# the new row with the new price-report entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("248:248").Insert()

$ws.Range("A248").Value = 4
$ws.Range("B248").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C248").Value = "Los Lagos"
$ws.Range("D248").Value = 44900
$ws.Range("E248").Value = 10
$ws.Range("F248").Value = 100112021
$ws.Range("G248").Value = "Ají"
$ws.Range("H248").Value = "Inferno"
$ws.Range("I248").Value = "Primera"
$ws.Range("J248").Value = 70
$ws.Range("K248").Value = 21000
$ws.Range("L248").Value = 21000
$ws.Range("M248").Value = 21000
$ws.Range("N248").Value = '$/caja 10 kilos'
$ws.Range("O248").Value = "Región de Arica y Parinacota"
$ws.Range("P248").Value = 2100
$ws.Range("Q248").Value = 10
$ws.Range("R248").Value = "Hortaliza"
